$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 / Row 4 / Row 7: new "test" cells (shared string index 5) ---
$ws.Range("D3").Value = "test"
$ws.Range("E4").Value = "test"
$ws.Range("E7").Value = "test"
$ws.Range("F7").Value = "test"

# Row 8: remove the formula cell that referenced D4 (=D4)
$ws.Range("D8").ClearContents()

# --- Row 9: plain variants with spaces ---
$ws.Range("A9").Value = "t ext"
$ws.Range("B9").Value = "te xt"
$ws.Range("C9").Value = "tex t"
$ws.Range("D9").Value = "t ex t"

# --- Row 10: newline variants (wrapped) ---
$ws.Range("A10").Value = "t`next"
$ws.Range("B10").Value = "te`nxt"
$ws.Range("C10").Value = "tex`nt"
$ws.Range("D10").Value = "text"

$ws.Range("A10:C10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 25.35

# --- Row 11: quote / comma variants ---
$ws.Range("A11").Value = "te`"xt"
$ws.Range("B11").Value = "te,xt"
$ws.Range("C11").Value = "te `" xt"
$ws.Range("D11").Value = "te , xt"

# --- Row 12: smart quotes / en dash / backslash variants ---
$ws.Range("A12").Value = [char]0x201C + "te-xt" + [char]0x201D
$ws.Range("B12").Value = "t e " + [char]0x2013 + " st"
$ws.Range("C12").Value = "te\st"
$ws.Range("D12").Value = "te\\st"

# --- View state: zoom + selection ---
$excel.ActiveWindow.Zoom = 280
$ws.Range("D10").Select()
